# Rename the first two sheets for the login test cases.
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

$ws1.Name = "ValidLogin"
$ws2.Name = "InvalidLogin"

# ValidLogin data set (headers + one valid admin/manager row).
$ws1.Range("A1").Value = "UserName"
$ws1.Range("B1").Value = "Password"
$ws1.Range("C1").Value = "eTitle"
$ws1.Range("A2").Value = "admin"
$ws1.Range("B2").Value = "manager"
$ws1.Range("C2").Value = "Enter"

# InvalidLogin data set (headers + one invalid abc/xyz row).
$ws2.Range("A1").Value = "UserName"
$ws2.Range("B1").Value = "Password"
$ws2.Range("A2").Value = "abc"
$ws2.Range("B2").Value = "xyz"

# Restore the selections recorded in each sheet view.
$ws1.Range("A1:C2").Select()
$ws2.Range("B3").Select()

# InvalidLogin is the active (selected) tab in the workbook.
$ws2.Activate()
